$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MASSO")

# Row 3 (Mar / Tuesday): replace the formula-driven "default availability" text
# with the specific masso schedule note (French / English).
$ws.Range("B3").Value = "Après la course - Claudie Germain"
$ws.Range("C3").Value = "After the race - Claudie Germain"

# Row 5 (Jeu / Thursday): re-point the formulas back to the base row (B2/C2)
# instead of chaining from row 4, restoring the default "according to
# availability" text for that day.
$ws.Range("B5").Formula = "=B2"
$ws.Range("C5").Formula = "=C2"
$ws.Range("C5").WrapText = $false
$ws.Rows.Item(5).AutoFit()

# Update the active selection to match the saved view state.
$ws.Range("B15").Select()
